$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Move the visitor-type dropdown values from column A (rows 3-11) to column B,
# leaving room for the new unique-id-type values.
for ($r = 3; $r -le 11; $r++) {
    $val = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $val
    $ws.Cells.Item($r, 1).ClearContents()
}

# New "Unique Id Type" supported values, listed under the existing
# "Unique Id Type" header in A17, in column B this time.
$ws.Cells.Item(18, 2).Value = "Adharcard"
$ws.Cells.Item(19, 2).Value = "Pancard"
$ws.Cells.Item(20, 2).Value = "Passport"

# Column B should match column A's width/bestFit formatting now that it
# holds the dropdown values too.
$colWidth = $ws.Columns.Item(1).ColumnWidth
$ws.Columns.Item(2).ColumnWidth = $colWidth

# Update the selection like Excel would leave it after this edit.
$ws.Range("E6").Select()
